# Insert a new data row at row 48 (sheet row), shifting existing rows 48-161 down to 49-162,
# and populate the new row 48 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 48; this shifts rows 48..161 down to 49..162
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows("48").Insert()

# Populate the newly inserted row 48 with the new record's values.
$ws.Cells.Item(48, 1).Value = 6
$ws.Cells.Item(48, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(48, 3).Value = "Metropolitana"
$ws.Cells.Item(48, 4).Value = 44622
$ws.Cells.Item(48, 5).Value = 13
$ws.Cells.Item(48, 6).Value = 100112029
$ws.Cells.Item(48, 7).Value = "Orégano"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 33
$ws.Cells.Item(48, 11).Value = 9000
$ws.Cells.Item(48, 12).Value = 10000
$ws.Cells.Item(48, 13).Value = 9455
$ws.Cells.Item(48, 14).Value = "$/docena de atados"
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value = 3152
$ws.Cells.Item(48, 17).Value = 3
$ws.Cells.Item(48, 18).Value = "Hortaliza"
